# Added item for order sets
#
# The "Items of the catalog" bullet list on slide 1 gets a new bullet —
# "Order sets: PlanDefinition" — inserted right after the
# "drug formulary: MedicationKnowledge," line and before the trailing
# "…" line. The textbox uses <a:spAutoFit/>, so PowerPoint grows its
# height automatically once the extra line is added; no manual resize
# is required.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "ZoneTexte 31" textbox (bulleted list of catalog item
# kinds) by name rather than a hard-coded index, so the script keeps
# working even if shape ordering ever shifts.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "ZoneTexte 31") {
        $targetShape = $candidate
        break
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph that currently reads "…" (it's the last bullet)
# and insert the new bullet immediately before it, so the new run
# inherits that paragraph's bullet/run formatting and becomes its own
# paragraph rather than merging into the previous run.
$lastParaIndex = $tr.Paragraphs().Count
$ellipsisPara = $tr.Paragraphs($lastParaIndex, 1)
[void]$ellipsisPara.InsertBefore("Order sets: PlanDefinition`r")
